$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "2026-02-14T17:00:31.471391+00:00"
$ws.Range("E2").Value = "Bhai bhej skte ho Thanks bro"
$ws.Range("F2").Value = "yes"
